$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 570.75
$ws.Range("J28").Value = 864.8
$ws.Range("L28").Value = 864.8
$ws.Range("N28").Value = -1834.8
$ws.Range("H40").Value = 1891.909
$ws.Range("I40").Value = 1928.8572
$ws.Range("K40").Value = 1928.8572
$ws.Range("M40").Value = -1753.8572
$ws.Range("H80").Value = 401.93024
$ws.Range("I80").Value = 282.75
$ws.Range("J80").Value = 505.56522
$ws.Range("K80").Value = 848.25
$ws.Range("L80").Value = 1516.69566
$ws.Range("M80").Value = 149.75
$ws.Range("N80").Value = -3512.69566
$ws.Range("H83").Value = 401.93024
$ws.Range("I83").Value = 282.75
$ws.Range("J83").Value = 505.56522
$ws.Range("K83").Value = 2544.75
$ws.Range("L83").Value = 4550.08698
$ws.Range("M83").Value = 2447.25
$ws.Range("N83").Value = -14534.08698
$ws.Range("H93").Value = 24800.5
$ws.Range("J93").Value = 24800.5
$ws.Range("L93").Value = 24800.5
$ws.Range("N93").Value = -29792.5
$ws.Range("H135").Value = 2558.838
$ws.Range("I135").Value = 2690.9
$ws.Range("J135").Value = 1992.8572
$ws.Range("K135").Value = 24218.1
$ws.Range("L135").Value = 17935.7148
$ws.Range("M135").Value = -21683.1
$ws.Range("N135").Value = -23005.7148
$ws.Range("H137").Value = 1673.279
$ws.Range("I137").Value = 1166.1
$ws.Range("J137").Value = 2114.3044
$ws.Range("K137").Value = 3498.3
$ws.Range("L137").Value = 6342.9132
$ws.Range("M137").Value = -948.2999999999997
$ws.Range("N137").Value = -11442.9132
$ws.Range("H138").Value = 3409.6545
$ws.Range("I138").Value = 1433.2963
$ws.Range("J138").Value = 5315.4287
$ws.Range("K138").Value = 4299.8889
$ws.Range("L138").Value = 15946.2861
$ws.Range("M138").Value = 840.1111000000001
$ws.Range("N138").Value = -26226.2861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1452.9615
$ws.Range("I2").Value = 1546.381
$ws.Range("J2").Value = 1060.6
$ws.Range("K2").Value = 1546.381
$ws.Range("L2").Value = 1060.6
$ws.Range("M2").Value = -1433.381
$ws.Range("N2").Value = -1286.6
$ws.Range("H32").Value = 4170.561
$ws.Range("I32").Value = 3720.3252
$ws.Range("J32").Value = 6661.8667
$ws.Range("K32").Value = 3720.3252
$ws.Range("L32").Value = 6661.8667
$ws.Range("M32").Value = -3433.3252
$ws.Range("N32").Value = -7235.8667
$ws.Range("H45").Value = 4594.607
$ws.Range("I45").Value = 8357
$ws.Range("J45").Value = 1333.8667
$ws.Range("K45").Value = 8357
$ws.Range("L45").Value = 1333.8667
$ws.Range("M45").Value = -7980
$ws.Range("N45").Value = -2087.8667
$ws.Range("H116").Value = 1452.9615
$ws.Range("I116").Value = 1546.381
$ws.Range("J116").Value = 1060.6
$ws.Range("K116").Value = 1546.381
$ws.Range("L116").Value = 1060.6
$ws.Range("M116").Value = 747.6189999999999
$ws.Range("N116").Value = -5648.6
$ws.Range("H139").Value = 45143
$ws.Range("J139").Value = 45143
$ws.Range("L139").Value = 45143
$ws.Range("N139").Value = -55423

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1452.9615
$ws.Range("I3").Value = 1546.381
$ws.Range("J3").Value = 1060.6
$ws.Range("K3").Value = 1546.381
$ws.Range("L3").Value = 1060.6
$ws.Range("M3").Value = -1432.381
$ws.Range("N3").Value = -1288.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2322.024
$ws.Range("I31").Value = 1394.3478
$ws.Range("J31").Value = 2671.8032
$ws.Range("K31").Value = 1394.3478
$ws.Range("L31").Value = 2671.8032
$ws.Range("M31").Value = -1099.3478
$ws.Range("N31").Value = -3261.8032
$ws.Range("H34").Value = 2322.024
$ws.Range("I34").Value = 1394.3478
$ws.Range("J34").Value = 2671.8032
$ws.Range("K34").Value = 1394.3478
$ws.Range("L34").Value = 2671.8032
$ws.Range("M34").Value = -1192.3478
$ws.Range("N34").Value = -3075.8032
$ws.Range("H93").Value = 19087.285
$ws.Range("I93").Value = 2722.2
$ws.Range("J93").Value = 60000
$ws.Range("K93").Value = 2722.2
$ws.Range("L93").Value = 60000
$ws.Range("M93").Value = -850.1999999999998
$ws.Range("N93").Value = -63744
$ws.Range("H99").Value = 15629602
$ws.Range("I99").Value = 1162.4
$ws.Range("J99").Value = 41677000
$ws.Range("K99").Value = 1162.4
$ws.Range("L99").Value = 41677000
$ws.Range("M99").Value = 335.5999999999999
$ws.Range("N99").Value = -41679996
$ws.Range("H122").Value = 1287.1177
$ws.Range("I122").Value = 1318.7
$ws.Range("J122").Value = 1242
$ws.Range("K122").Value = 3956.1
$ws.Range("L122").Value = 3726
$ws.Range("M122").Value = -1506.1
$ws.Range("N122").Value = -8626
$ws.Range("H126").Value = 15629602
$ws.Range("I126").Value = 1162.4
$ws.Range("J126").Value = 41677000
$ws.Range("K126").Value = 3487.2
$ws.Range("L126").Value = 125031000
$ws.Range("M126").Value = -1017.2
$ws.Range("N126").Value = -125035940

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 136388.56
$ws.Range("J5").Value = 163093.45
$ws.Range("L5").Value = 489280.35
$ws.Range("N5").Value = -489504.35
$ws.Range("H130").Value = 8275.454
$ws.Range("I130").Value = 3015
$ws.Range("J130").Value = 9444.444
$ws.Range("K130").Value = 9045
$ws.Range("L130").Value = 28333.332
$ws.Range("M130").Value = -4025
$ws.Range("N130").Value = -38373.33199999999
$ws.Range("H131").Value = 13096426
$ws.Range("J131").Value = 14926737
$ws.Range("L131").Value = 44780211
$ws.Range("N131").Value = -44790291
$ws.Range("H133").Value = 41692.793
$ws.Range("I133").Value = 127511.375
$ws.Range("K133").Value = 382534.125
$ws.Range("M133").Value = -377474.125
$ws.Range("H134").Value = 10452.767
$ws.Range("I134").Value = 11234.818
$ws.Range("K134").Value = 33704.454
$ws.Range("M134").Value = -28634.454
$ws.Range("H135").Value = 136388.56
$ws.Range("J135").Value = 163093.45
$ws.Range("L135").Value = 1467841.05
$ws.Range("N135").Value = -1472911.05
$ws.Range("H137").Value = 27792370
$ws.Range("I137").Value = 2461.7
$ws.Range("J137").Value = 38480796
$ws.Range("K137").Value = 7385.099999999999
$ws.Range("L137").Value = 115442388
$ws.Range("M137").Value = -2285.099999999999
$ws.Range("N137").Value = -115452588
$ws.Range("H139").Value = 4871.971
$ws.Range("I139").Value = 7304.3125
$ws.Range("J139").Value = 2823.6843
$ws.Range("K139").Value = 21912.9375
$ws.Range("L139").Value = 8471.052899999999
$ws.Range("M139").Value = -16772.9375
$ws.Range("N139").Value = -18751.0529
$ws.Range("H141").Value = 16847.842
$ws.Range("I141").Value = 13210.9
$ws.Range("J141").Value = 20888.889
$ws.Range("K141").Value = 39632.7
$ws.Range("L141").Value = 62666.667
$ws.Range("M141").Value = -34452.7
$ws.Range("N141").Value = -73026.667

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 58266.668
$ws.Range("J104").Value = 58266.668
$ws.Range("L104").Value = 58266.668
$ws.Range("N104").Value = -65254.668
$ws.Range("H122").Value = 5403734.5
$ws.Range("I122").Value = 5894255.5
$ws.Range("K122").Value = 17682766.5
$ws.Range("M122").Value = -17680316.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null

Write-Host "Applied all changes"